# Commit: "flipped N7 & S5 barcodes; fixed index2well bug"
#
# 1) Rename the two index-sheets so the tab order reads "N7-S5" instead of
#    "S5-N7" (the 515FB sheet is untouched).
# 2) On both the HiSeq and the MiSeq sheet, the "direction"/"well" columns
#    (D/E) had the forward (S5 "N7") and reverse (N7 "S5") blocks swapped,
#    and the per-block well numbering (F1/F2/F3 vs R1/R2) recomputed to fix
#    an index2well bug. Row order and columns A-C (sequence/index/name) are
#    unchanged; only D (direction) and E (well) move.

$wb = $excel.ActiveWorkbook

$wsHiSeq = $wb.Worksheets.Item("S5-N7_HiSeq")
$wsHiSeq.Name = "N7-S5_HiSeq"

$wsMiSeq = $wb.Worksheets.Item("S5-N7_MiSeq")
$wsMiSeq.Name = "N7-S5_MiSeq"

foreach ($ws in @($wsHiSeq, $wsMiSeq)) {
    $ws.Range("D2:D9").Value   = "forward"
    $ws.Range("E2:E9").Value   = "F1"

    $ws.Range("D10:D17").Value = "forward"
    $ws.Range("E10:E17").Value = "F2"

    $ws.Range("D18:D25").Value = "forward"
    $ws.Range("E18:E25").Value = "F3"

    $ws.Range("D26:D37").Value = "reverse"
    $ws.Range("E26:E37").Value = "R1"

    $ws.Range("D38:D49").Value = "reverse"
    $ws.Range("E38:E49").Value = "R2"
}

# Restore the per-sheet selections recorded in the saved view state.
# MiSeq is set first so the HiSeq tab ends up the active one (tabSelected).
$wsMiSeq.Range("D2").Select() | Out-Null
$wsHiSeq.Range("C38").Select() | Out-Null
